# "update for latest workshop sandisk"
#
# 1. The cached text of the auto-updating "date" fields (handout master,
#    notes master, slide master, and the slide layouts that carry their own
#    date placeholder) is refreshed from 6/18/2023 -> 6/22/2023.
# 2. The illustrative headshot picture that used to sit on the
#    "Key terminology" slide (slide 2) is removed.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, $oldText, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq $oldText) {
                $shape.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$oldDate = "6/18/2023"
$newDate = "6/22/2023"

# Handout master date placeholder
Set-DatePlaceholderText $p.HandoutMaster.Shapes $oldDate $newDate

# Notes master date placeholder
Set-DatePlaceholderText $p.NotesMaster.Shapes $oldDate $newDate

# Slide master date placeholder
Set-DatePlaceholderText $p.SlideMaster.Shapes $oldDate $newDate

# Slide layouts that carry their own cached date placeholder text
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $oldDate $newDate
}

# Remove the headshot picture from slide 2 ("Key terminology")
$slide = $p.Slides.Item(2)
for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Type -eq 13 -or $shape.Name -eq "Picture 4") {
        $shape.Delete()
    }
}
